$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 868.0952
$ws.Range("I15").Value = 868.0952
$ws.Range("K15").Value = 2604.2856
$ws.Range("M15").Value = -2435.2856
$ws.Range("H32").Value = 8867
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 8867
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 8867
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -9519
$ws.Range("H40").Value = 45461732
$ws.Range("I40").Value = 6713.6665
$ws.Range("K40").Value = 6713.6665
$ws.Range("M40").Value = -6538.6665
$ws.Range("H43").Value = 6098.2
$ws.Range("I43").Value = 6856.2856
$ws.Range("J43").Value = 4329.3335
$ws.Range("K43").Value = 6856.2856
$ws.Range("L43").Value = 4329.3335
$ws.Range("M43").Value = -6787.2856
$ws.Range("N43").Value = -4467.3335
$ws.Range("H51").Value = 6752.9414
$ws.Range("I51").Value = 8644.223
$ws.Range("K51").Value = 8644.223
$ws.Range("M51").Value = -8160.223
$ws.Range("H55").Value = 248.5
$ws.Range("I55").Value = 149.83333
$ws.Range("J55").Value = 347.16666
$ws.Range("K55").Value = 149.83333
$ws.Range("L55").Value = 347.16666
$ws.Range("M55").Value = 64.16667000000001
$ws.Range("N55").Value = -775.16666
$ws.Range("H96").Value = 1613355.6
$ws.Range("I96").Value = 2249.5
$ws.Range("J96").Value = 2073671.8
$ws.Range("K96").Value = 6748.5
$ws.Range("L96").Value = 6221015.4
$ws.Range("M96").Value = -5375.5
$ws.Range("N96").Value = -6223761.4
$ws.Range("H100").Value = 9147.091
$ws.Range("J100").Value = 12799.6
$ws.Range("L100").Value = 12799.6
$ws.Range("N100").Value = -13881.6
$ws.Range("H138").Value = 2908.0862
$ws.Range("J138").Value = 3149.8918
$ws.Range("L138").Value = 9449.6754
$ws.Range("N138").Value = -19729.6754

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 516.82355
$ws.Range("I2").Value = 449
$ws.Range("K2").Value = 449
$ws.Range("M2").Value = -336
$ws.Range("H5").Value = 737.5
$ws.Range("I5").Value = 700
$ws.Range("K5").Value = 700
$ws.Range("M5").Value = -588
$ws.Range("H32").Value = 6289.627
$ws.Range("I32").Value = 5667.6606
$ws.Range("K32").Value = 5667.6606
$ws.Range("M32").Value = -5380.6606
$ws.Range("H63").Value = 3616.3333
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 3616.3333
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H102").Value = 9704
$ws.Range("I102").Value = 9499
$ws.Range("K102").Value = 9499
$ws.Range("M102").Value = -7877
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H110").Value = 5098.7
$ws.Range("I110").Value = 4585.2666
$ws.Range("J110").Value = 6639
$ws.Range("K110").Value = 4585.2666
$ws.Range("L110").Value = 6639
$ws.Range("M110").Value = -2540.2666
$ws.Range("N110").Value = -10729
$ws.Range("H116").Value = 516.82355
$ws.Range("I116").Value = 449
$ws.Range("K116").Value = 449
$ws.Range("M116").Value = 1845

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 516.82355
$ws.Range("I3").Value = 449
$ws.Range("K3").Value = 449
$ws.Range("M3").Value = -335
$ws.Range("H4").Value = 737.5
$ws.Range("I4").Value = 700
$ws.Range("K4").Value = 700
$ws.Range("M4").Value = -585
$ws.Range("H64").Value = 993.3333
$ws.Range("I64").Value = 1000
$ws.Range("K64").Value = 1000
$ws.Range("M64").Value = -775
$ws.Range("H67").Value = 993.3333
$ws.Range("I67").Value = 1000
$ws.Range("K67").Value = 1000
$ws.Range("M67").Value = -220
$ws.Range("H105").Value = 432367.06
$ws.Range("I105").Value = 614372.7
$ws.Range("J105").Value = 7687.222
$ws.Range("K105").Value = 614372.7
$ws.Range("L105").Value = 7687.222
$ws.Range("M105").Value = -612625.7
$ws.Range("N105").Value = -11181.222

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1340.125
$ws.Range("I16").Value = 803.6667
$ws.Range("K16").Value = 803.6667
$ws.Range("M16").Value = -516.6667
$ws.Range("H105").Value = 7223.6665
$ws.Range("J105").Value = 15841.6
$ws.Range("L105").Value = 15841.6
$ws.Range("N105").Value = -19335.6
$ws.Range("H107").Value = 1154.7391
$ws.Range("I107").Value = 288.5
$ws.Range("J107").Value = 2502.2222
$ws.Range("K107").Value = 288.5
$ws.Range("L107").Value = 2502.2222
$ws.Range("M107").Value = 1631.5
$ws.Range("N107").Value = -6342.2222
$ws.Range("H113").Value = 1340.125
$ws.Range("I113").Value = 803.6667
$ws.Range("K113").Value = 803.6667
$ws.Range("M113").Value = 1366.3333

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 50000052
$ws.Range("I6").Value = 50000052
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 150000156
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -150000043
$ws.Range("N6").ClearContents()
$ws.Range("H82").Value = 11409.2
$ws.Range("I82").Value = 3499.5
$ws.Range("J82").Value = 16682.334
$ws.Range("K82").Value = 10498.5
$ws.Range("L82").Value = 50047.00199999999
$ws.Range("M82").Value = -10092.5
$ws.Range("N82").Value = -50859.00199999999
$ws.Range("H85").Value = 11409.2
$ws.Range("I85").Value = 3499.5
$ws.Range("J85").Value = 16682.334
$ws.Range("K85").Value = 10498.5
$ws.Range("L85").Value = 50047.00199999999
$ws.Range("M85").Value = -9094.5
$ws.Range("N85").Value = -52855.00199999999
$ws.Range("H132").Value = 1064.1111
$ws.Range("I132").Value = 796.7143
$ws.Range("K132").Value = 7170.428699999999
$ws.Range("M132").Value = -4640.428699999999
$ws.Range("H136").Value = 33333
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 33333
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 99999
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -110199

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1655.7407
$ws.Range("I102").Value = 1661.7307
$ws.Range("K102").Value = 1661.7307
$ws.Range("M102").Value = -39.73070000000007
$ws.Range("H107").Value = 786.3333
$ws.Range("I107").Value = 179.75
$ws.Range("K107").Value = 179.75
$ws.Range("M107").Value = 1740.25
$ws.Range("H118").Value = 45000
$ws.Range("J118").Value = 45000
$ws.Range("L118").Value = 45000
$ws.Range("N118").Value = -48314

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 866.38464
$ws.Range("I46").Value = 496
$ws.Range("K46").Value = 496
$ws.Range("M46").Value = -308
$ws.Range("H64").Value = 52426.6
$ws.Range("I64").Value = 20136
$ws.Range("K64").Value = 20136
$ws.Range("M64").Value = -19911
$ws.Range("H67").Value = 52426.6
$ws.Range("I67").Value = 20136
$ws.Range("K67").Value = 20136
$ws.Range("M67").Value = -19356
$ws.Range("H68").Value = 3476240.5
$ws.Range("I68").Value = 8334657
$ws.Range("J68").Value = 5943.143
$ws.Range("K68").Value = 8334657
$ws.Range("L68").Value = 5943.143
$ws.Range("M68").Value = -8333908
$ws.Range("N68").Value = -7441.143
$ws.Range("H71").Value = 3476240.5
$ws.Range("I71").Value = 8334657
$ws.Range("J71").Value = 5943.143
$ws.Range("K71").Value = 41673285
$ws.Range("L71").Value = 29715.715
$ws.Range("M71").Value = -41669541
$ws.Range("N71").Value = -37203.715

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 33199.6
$ws.Range("J2").Value = 16000
$ws.Range("L2").Value = 16000
$ws.Range("N2").Value = -16224
$ws.Range("H27").Value = 99999
$ws.Range("J27").Value = 99999
$ws.Range("L27").Value = 99999
$ws.Range("N27").Value = -100137
$ws.Range("H46").Value = 104147
$ws.Range("J46").Value = 104147
$ws.Range("L46").Value = 104147
$ws.Range("N46").Value = -104609
$ws.Range("H62").Value = 10471.667
$ws.Range("I62").Value = 3793.3333
$ws.Range("J62").Value = 17150
$ws.Range("K62").Value = 3793.3333
$ws.Range("L62").Value = 17150
$ws.Range("M62").Value = -3169.3333
$ws.Range("N62").Value = -18398
$ws.Range("H65").Value = 10471.667
$ws.Range("I65").Value = 3793.3333
$ws.Range("J65").Value = 17150
$ws.Range("K65").Value = 18966.6665
$ws.Range("L65").Value = 85750
$ws.Range("M65").Value = -15846.6665
$ws.Range("N65").Value = -91990
$ws.Range("H100").Value = 998
$ws.Range("I100").Value = 979.63635
$ws.Range("K100").Value = 1959.2727
$ws.Range("M100").Value = -1418.2727
$ws.Range("H109").Value = 120449.5
$ws.Range("J109").Value = 120449.5
$ws.Range("L109").Value = 120449.5
$ws.Range("N109").Value = -123223.5
$ws.Range("H115").Value = 92895.836
$ws.Range("J115").Value = 92895.836
$ws.Range("L115").Value = 92895.836
$ws.Range("N115").Value = -96029.836
$ws.Range("H132").Value = 403477.62
$ws.Range("I132").Value = 3542.7727
$ws.Range("J132").Value = 3336333.2
$ws.Range("K132").Value = 10628.3181
$ws.Range("L132").Value = 10008999.6
$ws.Range("M132").Value = -8098.3181
$ws.Range("N132").Value = -10014059.6
$ws.Range("H134").Value = 104147
$ws.Range("J134").Value = 104147
$ws.Range("L134").Value = 312441
$ws.Range("N134").Value = -317511
